# Generate Report for Handoff
# Updates the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" values
# for the bac28f38-e2d8-499d-92b0-c819b0177acc.md file, reflecting a newly
# generated handoff xliff report.

$wb = $excel.ActiveWorkbook

# Overview sheet: row 7 corresponds to bac28f38-e2d8-499d-92b0-c819b0177acc.md
# Column G = "Latest HO Xliff Generate Date"
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G7").Value = "2016-08-12 16:51:53"

# zh-cn sheet: row 7 corresponds to bac28f38-e2d8-499d-92b0-c819b0177acc.md
# Column H = "Latest Handoff Datetime"
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H7").Value = "2016-08-12 16:51:45"

# de-de sheet: row 7 corresponds to bac28f38-e2d8-499d-92b0-c819b0177acc.md
# Column H = "Latest Handoff Datetime"
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H7").Value = "2016-08-12 16:51:53"
